# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" for every cell that
#    carries that status (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2. Narrow the "Status" columns (Overview E:F, zh-cn C, de-de C) to match
#    the new, shorter status text.

$wb = $excel.ActiveWorkbook

$targetWidth = 13.4101845877511

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = $targetWidth
$overview.Columns.Item(6).ColumnWidth = $targetWidth

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = $targetWidth

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = $targetWidth
